$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts all existing
# columns (A:AS) one column to the right (B:AT), matching the target diff.
$ws.Columns.Item(1).Insert()

# Populate the new column A (rows 2-20) with a sequential 0-based row
# index, formatted the same way as the header row (bold font, thin box
# border, centered/top aligned) - i.e. the same direct formatting as the
# cellXfs style used by the header cells.
$rng = $ws.Range("A2:A20")
$rng.Borders.LineStyle = 1
$rng.Font.Bold = $true
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160

for ($i = 0; $i -lt 19; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}
